$wb = $excel.ActiveWorkbook

# ===== Sheet: LP1912 =====
$ws = $wb.Worksheets.Item("LP1912")
$ws.Cells.Item(2,1).Value = "Última actualización: 11:17:16"
$ws.Cells.Item(3,1).Value = "Total filas: 163"
$ws.Cells.Item(69,3).Value = "14_ABASTO"
$ws.Cells.Item(70,3).Value = "215A_EL PATO"
$ws.Cells.Item(73,1).Value = "07:48:14"
$ws.Cells.Item(73,3).Value = "215C_EL PATO"
$ws.Cells.Item(73,4).Value = 46
$ws.Cells.Item(74,1).Value = "08:30:59"
$ws.Cells.Item(74,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(74,4).Value = 4
$ws.Cells.Item(126,1).Value = "08:56:14"
$ws.Cells.Item(126,3).Value = "14_ABASTO"
$ws.Cells.Item(126,4).Value = 97
$ws.Cells.Item(127,1).Value = "10:26:25"
$ws.Cells.Item(127,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(127,4).Value = 7
$ws.Cells.Item(142,1).Value = "11:17:16"
$ws.Cells.Item(142,2).Value = "11:18"
$ws.Cells.Item(142,3).Value = "17_ROMERO"
$ws.Cells.Item(142,4).Value = 1
$ws.Cells.Item(143,2).Value = "11:19"
$ws.Cells.Item(143,3).Value = "215C_EL PATO"
$ws.Cells.Item(143,4).Value = 24
$ws.Cells.Item(144,1).Value = "10:55:35"
$ws.Cells.Item(144,2).Value = "11:20"
$ws.Cells.Item(144,4).Value = 25
$ws.Cells.Item(145,1).Value = "11:17:16"
$ws.Cells.Item(145,2).Value = "11:21"
$ws.Cells.Item(145,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(145,4).Value = 4
$ws.Cells.Item(146,1).Value = "11:17:16"
$ws.Cells.Item(146,4).Value = 16
$ws.Cells.Item(147,1).Value = "11:17:16"
$ws.Cells.Item(147,2).Value = "11:33"
$ws.Cells.Item(147,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(147,4).Value = 16
$ws.Cells.Item(148,1).Value = "11:17:16"
$ws.Cells.Item(148,2).Value = "11:41"
$ws.Cells.Item(148,3).Value = "16_SANTA ANA"
$ws.Cells.Item(148,4).Value = 24
$ws.Cells.Item(149,2).Value = "11:44"
$ws.Cells.Item(149,3).Value = "215B_EL PATO"
$ws.Cells.Item(149,4).Value = 49
$ws.Cells.Item(150,1).Value = "11:17:16"
$ws.Cells.Item(150,2).Value = "11:45"
$ws.Cells.Item(150,3).Value = "215B_EL PATO"
$ws.Cells.Item(150,4).Value = 28
$ws.Cells.Item(151,1).Value = "11:17:16"
$ws.Cells.Item(151,2).Value = "11:49"
$ws.Cells.Item(151,3).Value = "15_ABASTO"
$ws.Cells.Item(151,4).Value = 32
$ws.Cells.Item(152,1).Value = "11:17:16"
$ws.Cells.Item(152,2).Value = "11:51"
$ws.Cells.Item(152,3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(152,4).Value = 34
$ws.Cells.Item(153,1).Value = "11:17:16"
$ws.Cells.Item(153,2).Value = "11:56"
$ws.Cells.Item(153,3).Value = "225_GOMEZ"
$ws.Cells.Item(153,4).Value = 39
$ws.Cells.Item(154,1).Value = "11:17:16"
$ws.Cells.Item(154,2).Value = "12:01"
$ws.Cells.Item(154,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(154,4).Value = 44
$ws.Cells.Item(155,1).Value = "11:17:16"
$ws.Cells.Item(155,2).Value = "12:01"
$ws.Cells.Item(155,3).Value = "16_SANTA ANA"
$ws.Cells.Item(155,4).Value = 44
$ws.Cells.Item(156,1).Value = "11:17:16"
$ws.Cells.Item(156,2).Value = "12:04"
$ws.Cells.Item(156,3).Value = "17_ROMERO"
$ws.Cells.Item(156,4).Value = 47
$ws.Cells.Item(157,2).Value = "12:08"
$ws.Cells.Item(157,3).Value = "14_ABASTO"
$ws.Cells.Item(157,4).Value = 73
$ws.Cells.Item(158,1).Value = "11:17:16"
$ws.Cells.Item(158,2).Value = "12:09"
$ws.Cells.Item(158,3).Value = "14_ABASTO"
$ws.Cells.Item(158,4).Value = 52
$ws.Cells.Item(159,1).Value = "11:17:16"
$ws.Cells.Item(159,2).Value = "12:19"
$ws.Cells.Item(159,3).Value = "15_ABASTO"
$ws.Cells.Item(159,4).Value = 62
$ws.Cells.Item(160,1).Value = "11:17:16"
$ws.Cells.Item(160,2).Value = "12:20"
$ws.Cells.Item(160,3).Value = "10_OLMOS"
$ws.Cells.Item(160,4).Value = 63
$ws.Cells.Item(161,1).Value = "10:55:35"
$ws.Cells.Item(161,2).Value = "12:32"
$ws.Cells.Item(161,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(161,4).Value = 97
$ws.Cells.Item(161,5).Value = "LP1912"
$ws.Cells.Item(162,1).Value = "11:17:16"
$ws.Cells.Item(162,2).Value = "12:33"
$ws.Cells.Item(162,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(162,4).Value = 76
$ws.Cells.Item(162,5).Value = "LP1912"
$ws.Cells.Item(163,1).Value = "11:17:16"
$ws.Cells.Item(163,2).Value = "12:34"
$ws.Cells.Item(163,3).Value = "215C_EL PATO"
$ws.Cells.Item(163,4).Value = 77
$ws.Cells.Item(163,5).Value = "LP1912"
$ws.Cells.Item(164,1).Value = "11:17:16"
$ws.Cells.Item(164,2).Value = "12:36"
$ws.Cells.Item(164,3).Value = "27_EL RETIRO"
$ws.Cells.Item(164,4).Value = 79
$ws.Cells.Item(164,5).Value = "LP1912"
$ws.Cells.Item(165,1).Value = "11:17:16"
$ws.Cells.Item(165,2).Value = "12:47"
$ws.Cells.Item(165,3).Value = "10_OLMOS"
$ws.Cells.Item(165,4).Value = 90
$ws.Cells.Item(165,5).Value = "LP1912"
$ws.Cells.Item(166,1).Value = "11:17:16"
$ws.Cells.Item(166,2).Value = "12:51"
$ws.Cells.Item(166,3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(166,4).Value = 94
$ws.Cells.Item(166,5).Value = "LP1912"
$ws.Cells.Item(167,1).Value = "11:17:16"
$ws.Cells.Item(167,2).Value = "13:00"
$ws.Cells.Item(167,3).Value = "14_ABASTO"
$ws.Cells.Item(167,4).Value = 103
$ws.Cells.Item(167,5).Value = "LP1912"
$ws.Cells.Item(168,1).Value = "11:17:16"
$ws.Cells.Item(168,2).Value = "13:04"
$ws.Cells.Item(168,3).Value = "215C_EL PATO"
$ws.Cells.Item(168,4).Value = 107
$ws.Cells.Item(168,5).Value = "LP1912"

# ===== Sheet: LP1912-215 =====
$ws = $wb.Worksheets.Item("LP1912-215")
$ws.Cells.Item(2,1).Value = "Última actualización: 11:17:16"
$ws.Cells.Item(3,1).Value = "Total filas: 31"
$ws.Cells.Item(23,3).Value = "215B_EL PATO"
$ws.Cells.Item(24,3).Value = "215A_EL PATO"
$ws.Cells.Item(32,1).Value = "11:17:16"
$ws.Cells.Item(32,4).Value = 16
$ws.Cells.Item(34,1).Value = "11:17:16"
$ws.Cells.Item(34,2).Value = "11:45"
$ws.Cells.Item(34,3).Value = "215B_EL PATO"
$ws.Cells.Item(34,4).Value = 28
$ws.Cells.Item(35,1).Value = "11:17:16"
$ws.Cells.Item(35,2).Value = "12:34"
$ws.Cells.Item(35,3).Value = "215C_EL PATO"
$ws.Cells.Item(35,4).Value = 77
$ws.Cells.Item(35,5).Value = "LP1912"
$ws.Cells.Item(36,1).Value = "11:17:16"
$ws.Cells.Item(36,2).Value = "13:04"
$ws.Cells.Item(36,3).Value = "215C_EL PATO"
$ws.Cells.Item(36,4).Value = 107
$ws.Cells.Item(36,5).Value = "LP1912"

# ===== Sheet: 6203-6173 =====
$ws = $wb.Worksheets.Item("6203-6173")
$ws.Cells.Item(2,1).Value = "Última actualización: 11:17:16"
$ws.Cells.Item(3,1).Value = "Total filas: 28"
$ws.Cells.Item(31,1).Value = "11:17:16"
$ws.Cells.Item(31,4).Value = 39
$ws.Cells.Item(32,1).Value = "11:17:16"
$ws.Cells.Item(32,2).Value = "12:56"
$ws.Cells.Item(32,3).Value = "215C_LA PLATA"
$ws.Cells.Item(32,4).Value = 99
$ws.Cells.Item(32,5).Value = "L6203"
$ws.Cells.Item(33,1).Value = "11:17:16"
$ws.Cells.Item(33,2).Value = "13:11"
$ws.Cells.Item(33,3).Value = "215A_LA PLATA"
$ws.Cells.Item(33,4).Value = 114
$ws.Cells.Item(33,5).Value = "L6173"
